$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Insert two blank columns before old column M ("Huong") to make room for
#    the new "So do tang day" / "Phoi canh" columns (L = "Hinh anh" keeps its
#    position and gets repurposed into "Bo suu tap").
$ws.Range("M1:N1").EntireColumn.Insert()

# 2. Update header text for K, L, M, N.
$ws.Range("K1").Value = "Hình đại diện (1 hình duy nhất)"
$ws.Range("L1").Value = "Bộ sưu tập"
$ws.Range("M1").Value = "Sơ đồ tầng dãy"
$ws.Range("N1").Value = "Phối cảnh"

# 3. Populate the two new data cells in row 2 with the same sample value that
#    the original "Hinh anh" column (L2) already carries.
$imgList = $ws.Range("L2").Value
$ws.Range("M2").Value = $imgList
$ws.Range("N2").Value = $imgList

# 4. Header row formatting: bold header keeps its fill/vertical-center, but
#    now also wraps text; K1 gets its own (slightly darker) fill to stand out.
$ws.Range("A1:AB1").WrapText = $true
$ws.Range("K1").Interior.ThemeColor = 1
$ws.Range("K1").Interior.TintAndShade = -0.14999847407452621

# 5. Body row: wrap text on every populated cell.
$ws.Range("A2:AB2").WrapText = $true

# 6. Column L is much wider now that it holds free-text image lists.
$ws.Columns.Item(12).ColumnWidth = 62

# 7. Row 2 grows to fit the wrapped long-form text.
$ws.Rows.Item(2).RowHeight = 409.5

# 8. Reset the view: selection on L2.
$ws.Range("L2").Select()
